# Edit matching the commit diff:
#   - Slide 7's subtitle shape: "We've gone over a lot since the last quiz"
#     becomes two runs: "We've gone over a lot since the last " + "review"
#
# (The embedded-font typeface swap in the diff (Lato <-> Raleway metadata in
# <p:embeddedFontLst>) has no surface in the PowerPoint object model -- there
# is no ActivePresentation property/collection that exposes or edits the
# embedded font table's typeface names, so it cannot be produced via COM
# automation and is intentionally left alone here.)

$p = $ppt.ActivePresentation

# Locate the shape on slide 7 that holds the subtitle text we need to edit,
# rather than hard-coding a shape index.
$slide = $p.Slides.Item(7)
$targetShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text -like "*gone over a lot*") {
            $targetShape = $candidate
            break
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Replace only the trailing word "quiz" with "review", leaving the rest of
# the run (including the curly apostrophe in "We've") completely untouched.
# Using Characters() on the sub-range splits the paragraph into two runs at
# exactly the boundary the diff shows, without touching any other
# formatting.
$fullText = $tr.Text
$oldWord = "quiz"
$newWord = "review"
$wordStart = $fullText.LastIndexOf($oldWord) + 1

$wordRange = $tr.Characters($wordStart, $oldWord.Length)
$wordRange.Text = $newWord
